# Remove files and materials modification
# Adds a new "Raspberry PI 3 B+" line item (Farnell part 2842228) to the
# materials list, pushing the TOTAL row (and the trailing helper row) down
# by one, and updates the TOTAL formula to include the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right above the current TOTAL row (row 21), shifting
# TOTAL (row 21 -> 22) and the trailing helper row (row 35 -> 36) down.
$ws.Rows(21).Insert()

# Copy the formatting of the previous data row (20) onto the new row so the
# new cells get the same "Normal" styling used throughout the table.
$ws.Range("A20:E20").Copy()
$ws.Range("A21").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K20").Copy()
$ws.Range("K21").PasteSpecial(-4122)   # xlPasteFormats (keeps the spacer cell style)

# New line item data.
$ws.Range("A21").Value = "Raspberry PI 3 B+"
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 2842228
$ws.Range("D21").Value = "https://es.farnell.com/raspberry-pi/rpi3-modbp/ordenador-monoplaca-raspberry/dp/2842228"
$ws.Range("E21").Value = 32.46

# Turn D21 into a real hyperlink (matches the style used by D3/D20).
$ws.Hyperlinks.Add($ws.Range("D21"), "https://es.farnell.com/raspberry-pi/rpi3-modbp/ordenador-monoplaca-raspberry/dp/2842228")

# Re-apply the hyperlink-cell formatting so D21 reuses the existing
# "Hipervinculo" style (same as D20) instead of a freshly minted one.
$ws.Range("D20").Copy()
$ws.Range("D21").PasteSpecial(-4122)   # xlPasteFormats

# Extend the TOTAL formula to include the newly inserted row.
$ws.Range("E22").Formula = "=SUM(E2:E21)"

# Page setup metadata present in the target workbook.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Match the final cell selection recorded in the saved workbook.
$ws.Range("E22").Select() | Out-Null
